$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Rows 22 and 23 hold the same observation fields except for the
# identifying / species data -- this edit swaps that data between
# the two rows (row 22 becomes the old row 23 and vice versa).
# -----------------------------------------------------------------

# Save old values first so we can cross-assign them.
$A22 = $ws.Range("A22").Value2
$B22 = $ws.Range("B22").Value2
$D22 = $ws.Range("D22").Value2
$E22 = $ws.Range("E22").Value2
$F22 = $ws.Range("F22").Value2
$G22 = $ws.Range("G22").Value2
$H22 = $ws.Range("H22").Value2
$Q22 = $ws.Range("Q22").Value2
$R22 = $ws.Range("R22").Value2

$A23 = $ws.Range("A23").Value2
$B23 = $ws.Range("B23").Value2
$D23 = $ws.Range("D23").Value2
$E23 = $ws.Range("E23").Value2
$F23 = $ws.Range("F23").Value2
$G23 = $ws.Range("G23").Value2
$H23 = $ws.Range("H23").Value2
$Q23 = $ws.Range("Q23").Value2
$R23 = $ws.Range("R23").Value2

# Row 22 <- old row 23 data
$ws.Range("A22").Value = $A23
$ws.Range("B22").Value = $B23
$ws.Range("D22").Value = $D23
$ws.Range("E22").Value = $E23
$ws.Range("F22").Value = $F23
$ws.Range("G22").Value = $G23
$ws.Range("H22").Value = $H23
$ws.Range("I22").Value = "1"
$ws.Range("Q22").Value = $Q23
$ws.Range("R22").Value = $R23

# Row 23 <- old row 22 data
$ws.Range("A23").Value = $A22
$ws.Range("B23").Value = $B22
$ws.Range("D23").Value = $D22
$ws.Range("E23").Value = $E22
$ws.Range("F23").Value = $F22
$ws.Range("G23").Value = $G22
$ws.Range("H23").Value = $H22
$ws.Range("I23").ClearContents()
$ws.Range("Q23").Value = $Q22
$ws.Range("R23").Value = $R22

# -----------------------------------------------------------------
# Rows 26 and 27 likewise swap their species / location / date data.
# -----------------------------------------------------------------

$A26 = $ws.Range("A26").Value2
$B26 = $ws.Range("B26").Value2
$E26 = $ws.Range("E26").Value2
$F26 = $ws.Range("F26").Value2
$G26 = $ws.Range("G26").Value2
$Q26 = $ws.Range("Q26").Value2
$R26 = $ws.Range("R26").Value2
$Y26 = $ws.Range("Y26").Value2
$AA26 = $ws.Range("AA26").Value2

$A27 = $ws.Range("A27").Value2
$B27 = $ws.Range("B27").Value2
$E27 = $ws.Range("E27").Value2
$F27 = $ws.Range("F27").Value2
$G27 = $ws.Range("G27").Value2
$Q27 = $ws.Range("Q27").Value2
$R27 = $ws.Range("R27").Value2
$Y27 = $ws.Range("Y27").Value2
$AA27 = $ws.Range("AA27").Value2

# Row 26 <- old row 27 data
$ws.Range("A26").Value = $A27
$ws.Range("B26").Value = $B27
$ws.Range("E26").Value = $E27
$ws.Range("F26").Value = $F27
$ws.Range("G26").Value = $G27
$ws.Range("Q26").Value = $Q27
$ws.Range("R26").Value = $R27
# Dates are stored as plain text in this sheet; force text formatting
# so Excel doesn't reinterpret the "yyyy-mm-dd" string as a real date.
$ws.Range("Y26").NumberFormat = "@"
$ws.Range("Y26").Value = $Y27
$ws.Range("AA26").NumberFormat = "@"
$ws.Range("AA26").Value = $AA27

# Row 27 <- old row 26 data
$ws.Range("A27").Value = $A26
$ws.Range("B27").Value = $B26
$ws.Range("E27").Value = $E26
$ws.Range("F27").Value = $F26
$ws.Range("G27").Value = $G26
$ws.Range("Q27").Value = $Q26
$ws.Range("R27").Value = $R26
$ws.Range("Y27").NumberFormat = "@"
$ws.Range("Y27").Value = $Y26
$ws.Range("AA27").NumberFormat = "@"
$ws.Range("AA27").Value = $AA26
